$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix C28: was stored as text "04252502000160" (with leading zero),
# should become the numeric value 4252502000160 (matching the other CNPJ cells)
$ws.Range("C28").NumberFormat = "General"
$ws.Range("C28").Value = 4252502000160

# --- Append new log rows (29-32) ---

# Row 29
$ws.Range("A29").Value = "03/07/2025 21:40:09"
$ws.Range("B29").Value = "Ima Industria"
$ws.Range("C29").NumberFormat = "General"
$ws.Range("C29").Value = 4252502000160
$ws.Range("D29").Value = "denissonfhsilva@gmail.com"
$ws.Range("E29").Value = "893-ExtratoMensal-052025.pdf"

# Row 30
$ws.Range("A30").Value = "03/07/2025 22:57:37"
$ws.Range("B30").Value = "Ima Industria"
$ws.Range("C30").NumberFormat = "General"
$ws.Range("C30").Value = 4252502000160
$ws.Range("D30").Value = "denissonfhsilva@gmail.com"
$ws.Range("E30").Value = "893-ExtratoMensal-052025.pdf"

# Row 31
$ws.Range("A31").Value = "04/07/2025 11:20:05"
$ws.Range("B31").Value = "Ima Industria"
$ws.Range("C31").NumberFormat = "General"
$ws.Range("C31").Value = 4252502000160
$ws.Range("D31").Value = "denissonfhsilva@gmail.com"
$ws.Range("E31").Value = "893-ExtratoMensal-052025.pdf"

# Row 32 (note: CNPJ stays as text here, matching the source diff)
$ws.Range("A32").Value = "04/07/2025 11:20:15"
$ws.Range("B32").Value = "RCD Educação"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "58475425000137"
$ws.Range("D32").Value = "denissonfhsilva@gmail.com"
$ws.Range("E32").Value = "901-ExtratoMensal-052025.pdf"

Write-Output "done"
